$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rename the worksheet to reflect the survey's budget data source
$ws.Name = "FRB_Budget"

# Scroll the frozen view down and select the row the user is editing
[void]$ws.Range("A30").Select()
$excel.ActiveWindow.ScrollRow = 30

# Change the answer-type for question 49 from "number" to "radio"
$ws.Range("B49").Value = "radio"
[void]$ws.Range("B49").Select()
